$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 87
$ws1.Range("F3").Value = 11882
$ws1.Range("F4").Value = 18
$ws1.Range("F7").Value = 225
$ws1.Range("F8").Value = 11801
$ws1.Range("F9").Value = 488
$ws1.Range("F12").Value = 57
$ws1.Range("F13").Value = 1778
$ws1.Range("F14").Value = 5848
$ws1.Range("F15").Value = 124
$ws1.Range("F16").Value = 3536
$ws1.Range("F17").Value = 189
$ws1.Range("F18").Value = 21

# Sheet "演出" (performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 575

# Sheet "全部类型" (all types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 575
$ws4.Range("F3").Value = 87
$ws4.Range("F5").Value = 11882
$ws4.Range("F6").Value = 18
$ws4.Range("F10").Value = 225
$ws4.Range("F11").Value = 11801
$ws4.Range("F12").Value = 488
$ws4.Range("F15").Value = 57
$ws4.Range("F16").Value = 1778
$ws4.Range("F18").Value = 5848
$ws4.Range("F19").Value = 124
$ws4.Range("F20").Value = 3536
$ws4.Range("F21").Value = 189
$ws4.Range("F22").Value = 21
